# Update "想去人数" (F column) counts that changed for this data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2420
$ws.Range("F4").Value = 31
$ws.Range("F6").Value = 69
$ws.Range("F7").Value = 281
$ws.Range("F9").Value = 3271
$ws.Range("F10").Value = 1165
$ws.Range("F12").Value = 863
$ws.Range("F14").Value = 844
$ws.Range("F15").Value = 1492
$ws.Range("F16").Value = 766
$ws.Range("F19").Value = 373
$ws.Range("F20").Value = 71
$ws.Range("F21").Value = 114
$ws.Range("F23").Value = 2636

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F19").Value = 155
$ws.Range("F26").Value = 46
$ws.Range("F38").Value = 355

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2487
$ws.Range("F6").Value = 2504
$ws.Range("F7").Value = 9541
$ws.Range("F8").Value = 139
$ws.Range("F10").Value = 16
$ws.Range("F11").Value = 362
$ws.Range("F12").Value = 2803
$ws.Range("F13").Value = 368
$ws.Range("F14").Value = 680

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2487
$ws.Range("F3").Value = 139
$ws.Range("F6").Value = 2420
$ws.Range("F7").Value = 362
$ws.Range("F8").Value = 368
$ws.Range("F13").Value = 69
$ws.Range("F14").Value = 281
$ws.Range("F16").Value = 1165
$ws.Range("F19").Value = 863
$ws.Range("F21").Value = 844
$ws.Range("F23").Value = 1492
$ws.Range("F27").Value = 766
$ws.Range("F32").Value = 373
$ws.Range("F35").Value = 46
$ws.Range("F39").Value = 71
$ws.Range("F44").Value = 355
$ws.Range("F45").Value = 2636
